# Apply "changed some ROI parameters" edit to the BEC experiment-type table.
#
# Row 21 (NiLatticeSlosh) gets its RoiName switched from "Bec" to "NiLattice",
# its AnalysisMethod trimmed (drop "Tof"), and its CenterFitMethod switched
# from ParabolicFit1D to SineFit1D.
#
# Two brand-new rows are appended describing additional ROI-driven trial
# configurations: row 22 (NiLattice / RunIndex) and row 23 (NiBec / RunIndex).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: update existing trial (NiLatticeSlosh) ---
$ws.Range("D21").Value = "NiLattice"
$ws.Range("G21").Value = "DensityFit;AtomNumber;CenterFit"
$ws.Range("N21").Value = "SineFit1D"

# --- Row 22: new trial (NiLattice) ---
$ws.Range("A22").Value = "NiLattice"
$ws.Range("B22").Value = "A slosh experiment at the non-interacting lattice stage."
$ws.Range("C22").Value = "TOP"
$ws.Range("D22").Value = "NiLattice"
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = "RunIndex"
$ws.Range("G22").Value = "DensityFit;AtomNumber"
$ws.Range("H22").Value = "LSR"
$ws.Range("I22").Value = "HF"
$ws.Range("J22").Value = "StrongLight"
$ws.Range("K22").Value = 8
$ws.Range("L22").Value = "BosonicGaussianFit1D"
$ws.Range("M22").Value = 1
$ws.Range("N22").Value = "ParabolicFit1D"

# --- Row 23: new trial (NiBec) ---
$ws.Range("A23").Value = "NiBec"
$ws.Range("B23").Value = "A TOF experiment at the non-interacting BEC stage."
$ws.Range("C23").Value = "TOP"
$ws.Range("D23").Value = "Bec"
$ws.Range("E23").Value = 4
$ws.Range("F23").Value = "RunIndex"
$ws.Range("G23").Value = "CenterFit;AtomNumber;DensityFit"
$ws.Range("H23").Value = "LSR"
$ws.Range("I23").Value = "HF"
$ws.Range("J23").Value = "StrongLight"
$ws.Range("K23").Value = 8
$ws.Range("L23").Value = "BosonicGaussianFit1D"
$ws.Range("M23").Value = 1
$ws.Range("N23").Value = "ParabolicFit1D"
